$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column G header + existing rows
$ws.Range("G1").Value = "Assigner"
$ws.Range("G2").Value = "ds"
$ws.Range("G3").Value = "dhanasekar"
$ws.Range("G4").Value = "ds"

# New row 5 data
$ws.Range("A5").Value = "assigneng4"
$ws.Range("B5").Value = "dans1"
$ws.Range("C5").Value = "wheels india"
$ws.Range("D5").Value = 43132
$ws.Range("D2").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("E5").Value = "once a week"
$ws.Range("F5").Value = "Project"
$ws.Range("G5").Value = "ds"

# Move selection to H5, matching the post-edit cursor position
$ws.Range("H5").Select()
